# Plano Projeto.xlsx -- "Conexao com DB, Insercao de dados nas tabelas"
#
# The commit adds a new "OK" status column (C) for the plan's checklist rows,
# adds four new plan steps (7th - 10th) with their task names, and fills in
# the previously-empty 6th-step task name. It also moves the active
# selection further down the sheet (to B23) and scrolls the view so row 7
# is back at the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 6th step gains a task name (column B), plus "OK" status in column C
#     for every completed step (rows 16-22) ---
$ws.Range("B21").Value = "responsivo"

$ws.Range("C16").Value = "OK"
$ws.Range("C17").Value = "OK"
$ws.Range("C18").Value = "OK"
$ws.Range("C19").Value = "OK"
$ws.Range("C20").Value = "OK"
$ws.Range("C21").Value = "OK"
$ws.Range("C22").Value = "OK"

# --- New plan rows: 7th - 10th steps (columns A and B) ---
$ws.Range("A22").Value = "7º"
$ws.Range("B22").Value = "métodos PDO"

$ws.Range("A23").Value = "8º"
$ws.Range("B23").Value = "conf. Search"

$ws.Range("A24").Value = "9º"
$ws.Range("B24").Value = "Paginas"

$ws.Range("A25").Value = "10º"
$ws.Range("B25").Value = "Site dinâmico"

# --- View state: scroll back up a bit and land the selection on B23 ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B23").Select()
